$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): I0 in I1, IF in J1 -- same header formatting as the
# other header cells (e.g. H1: bold font, border, centered alignment).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# New data cells (row 2): numeric values 9 for I2 and J2, unstyled like H2.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
